$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New 'K' (strikeout) values regenerated for column G, keyed by worksheet row number.
$kValues = @{
    2 = 2
    3 = 2
    4 = 2
    5 = 1
    6 = 3
    7 = 1
    8 = 1
    9 = 1
    10 = 1
    11 = 0
    12 = 1
    13 = 0
    14 = 2
    15 = 1
    16 = 0
    17 = 1
    18 = 1
    19 = 1
    20 = 2
    21 = 2
    22 = 0
    23 = 0
    24 = 2
    25 = 1
    26 = 0
    27 = 2
    28 = 0
    29 = 1
    30 = 1
    31 = 1
    32 = 1
    33 = 2
    34 = 2
    35 = 2
    36 = 2
    37 = 0
    38 = 0
    39 = 2
    40 = 2
    41 = 2
    42 = 0
    43 = 2
    44 = 1
    45 = 0
    46 = 0
    47 = 2
    48 = 0
    49 = 1
    50 = 0
    51 = 1
    52 = 1
    53 = 0
    54 = 2
    55 = 0
    56 = 1
    57 = 1
    58 = 0
    59 = 1
    60 = 0
    61 = 0
    62 = 2
    63 = 2
    64 = 1
    65 = 0
    66 = 3
    67 = 2
    69 = 1
    70 = 1
    71 = 0
    72 = 1
}

foreach ($row in $kValues.Keys) {
    $ws.Cells.Item($row, 7).Value = $kValues[$row]
}
